$d = $word.ActiveDocument

# =======================================================================
# 1) "Bao gồm các page sau:" paragraph -> mark text as French (fr-FR).
#    (paragraph index 4; text content is unchanged)
# =======================================================================
$d.Paragraphs(4).Range.LanguageID = "fr-FR"

# =======================================================================
# 2) "Page Nhân Viên hiển thị tất cả nhân viên ..." paragraph
#    -> mark text as French (fr-FR). (paragraph index 5; text unchanged)
# =======================================================================
$d.Paragraphs(5).Range.LanguageID = "fr-FR"

# =======================================================================
# 3) "Chức năng hiển thị thông tin chi tiết ..." paragraph
#    -> mark text as French (fr-FR). (paragraph index 6; text unchanged)
# =======================================================================
$d.Paragraphs(6).Range.LanguageID = "fr-FR"

# =======================================================================
# 4) "Page Phòng ban hiển thị tất cả các phong ban ..." paragraph
#    -> fix missing diacritic: "phong" -> "phòng".
# =======================================================================
$d.Content.Find.Execute("hiển thị tất cả các phong ban", $true, $false, $false, $false, $false, $true, 1, $false, "hiển thị tất cả các phòng ban", 2) | Out-Null

# =======================================================================
# 5) "Gồm chức năng sắp xếp theo tên nhân viên." paragraph
#    -> change "tên" to "mã".
# =======================================================================
$d.Content.Find.Execute("Gồm chức năng sắp xếp theo tên nhân viên.", $true, $false, $false, $false, $false, $true, 1, $false, "Gồm chức năng sắp xếp theo mã nhân viên.", 2) | Out-Null

# =======================================================================
# 6) "Page Nhân Viên." paragraph (paragraph index 12)
#    -> append " (khi lựa chọn menu nhân viên)" and mark as French.
# =======================================================================
$r12 = $d.Paragraphs(12).Range
$ins12 = $d.Range($r12.End - 1, $r12.End - 1)
$ins12.InsertAfter(" (khi lựa chọn menu nhân viên)")
$d.Paragraphs(12).Range.LanguageID = "fr-FR"

# =======================================================================
# 7) "Page phòng ban" paragraph (paragraph index 16)
#    -> append " (hiển thị khi lựa chọn menu phòng ban)".
# =======================================================================
$r16 = $d.Paragraphs(16).Range
$ins16 = $d.Range($r16.End - 1, $r16.End - 1)
$ins16.InsertAfter(" (hiển thị khi lựa chọn menu phòng ban)")

# =======================================================================
# 8) "Page bảng lương" paragraph (paragraph index 18)
#    -> append " (hiển thị khi lựa chọn menu bảng lương)".
# =======================================================================
$r18 = $d.Paragraphs(18).Range
$ins18 = $d.Range($r18.End - 1, $r18.End - 1)
$ins18.InsertAfter(" (hiển thị khi lựa chọn menu bảng lương)")

Write-Host "Edits applied successfully."
